# Update TPM values for rows 2-4 and remove the "Resolving-Mac" row (row 5)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6855646666666667
$ws.Range("H2").Value = 2.056694
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.124002
$ws.Range("N2").Value = 0.372006
$ws.Range("O2").Value = 0.01161425268110074
$ws.Range("P2").Value = 0.01161425268110074
$ws.Range("Q2").Value = 0.08501138979600001
$ws.Range("R2").Value = 0.7651025081640001
$ws.Range("S2").Value = 0.01161425268110074
$ws.Range("T2").Value = 0.01161425268110074

# Row 3 (FAPs)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6855646666666667
$ws.Range("H3").Value = 2.056694
$ws.Range("O3").Value = 0.5198585843927942
$ws.Range("P3").Value = 0.5198585843927942
$ws.Range("Q3").Value = 3.805143729008666
$ws.Range("R3").Value = 34.246293561078
$ws.Range("S3").Value = 0.5198585843927942
$ws.Range("T3").Value = 0.5198585843927942

# Row 4 (MuSCs)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6855646666666667
$ws.Range("H4").Value = 2.056694
$ws.Range("M4").Value = 5.002328333333334
$ws.Range("N4").Value = 15.006985
$ws.Range("O4").Value = 0.4685271629261051
$ws.Range("P4").Value = 0.4685271629261051
$ws.Range("Q4").Value = 3.429419556398889
$ws.Range("R4").Value = 30.86477600759
$ws.Range("S4").Value = 0.4685271629261051
$ws.Range("T4").Value = 0.4685271629261051

# Remove row 5 (Target cluster = Resolving-Mac) entirely
$ws.Rows("5:5").Delete()
